$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 538, shifting existing rows 538:568 down to 539:569
$ws.Rows.Item(538).Insert()

# Populate the newly inserted row with its data
$ws.Cells.Item(538, 1).Value = 10
$ws.Cells.Item(538, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(538, 3).Value = 'La Araucanía'
$ws.Cells.Item(538, 4).Value = 45267
$ws.Cells.Item(538, 5).Value = 9
$ws.Cells.Item(538, 6).Value = 100114013
$ws.Cells.Item(538, 7).Value = 'Zanahoria'
$ws.Cells.Item(538, 8).Value = 'Sin especificar'
$ws.Cells.Item(538, 9).Value = 'Primera'
$ws.Cells.Item(538, 10).Value = 100
$ws.Cells.Item(538, 11).Value = 7000
$ws.Cells.Item(538, 12).Value = 7000
$ws.Cells.Item(538, 13).Value = 7000
$ws.Cells.Item(538, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(538, 15).Value = 'Región del Maule'
$ws.Cells.Item(538, 16).Value = 350
$ws.Cells.Item(538, 17).Value = 20
$ws.Cells.Item(538, 18).Value = 'Hortaliza'
